$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.459.79"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.727.09"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D5").Value = "'243.61"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.4870"
$ws.Range("D8").Value = "'0.2623"
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("D9").Value = "'0.06174"
$ws.Range("E9").Value = "  -0.63%  "
$ws.Range("D10").Value = "1.733.56"
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'0.07030"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "'15.47"
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").Value = "'4.568"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'0.5998"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").Value = "'77.28"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "26.481.14"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'0.000007075"
$ws.Range("E19").Value = "  +2.12%  "
$ws.Range("E20").Value = "  -1.92%  "
$ws.Range("D21").Value = "1.961.25"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'4.471"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "'8.603"
$ws.Range("E23").Value = "  -3.70%  "
$ws.Range("D24").Value = "'5.185"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'139.39"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "'15.29"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'1.406"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").Value = "'106.63"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "'1.713"
$ws.Range("E29").Value = "  -4.46%  "
$ws.Range("D30").Value = "'3.960"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").Value = "'0.07978"
$ws.Range("E31").Value = "  -0.54%  "
$ws.Range("D32").Value = "'3.687"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  -0.77%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.614"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.000"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6238"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.9094"
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").Value = "'1.984"
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.408"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "'1.002"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.01481"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'100.21"
$ws.Range("E42").Value = "  -4.72%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.436"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.3860"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'6.679"
$ws.Range("E45").Value = "  -3.73%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "'0.1157"
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05366"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "'30.28"
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.708"
$ws.Range("E49").Value = "  -1.69%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.251"
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'51.03"
$ws.Range("E51").Value = "  -0.39%  "
